# Generate Report for Handoff
#
# A new handoff was just kicked off for the "ed17b87a-b2e2-48ad-8870-0ee451091651"
# file, so its "Latest Handoff Datetime" needs to be refreshed on both the
# zh-cn and de-de localization-status sheets (column D, row 3 -- the row for
# that file). Everything else on the report stays the same.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-11 07:54:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-11 07:54:23"
